# fix Ef_En_PU_Prf3 func call bug in penv
# Correct the wr_cr_off values in column M (rows 2 and 3) which were
# computed incorrectly. A leading apostrophe is used so Excel stores the
# corrected figures as text (matching how the original, correct values
# were stored) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = "'-0.1006"
$ws.Range("M3").Value = "'-0.1504"
